# Applies the "adding averages and more checks" commit:
#  - Refresh the LAST UPDATE date (08-Sep-2025 -> 16-Sep-2025) and recompute
#    PERIOD TO EXPIRE (-8 days) on the Training Dashboard sheet.
#  - Row 17 (LOTO SOP) has now expired: flip it to the NOT VALID look/values.
#  - Restyle the title / header font (bold white) on both sheets.
#  - Update the Exam Dashboard COMMENTS column and narrow its column.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Training Dashboard
$ws2 = $wb.Worksheets.Item(2)   # Exam Dashboard

# ---------------------------------------------------------------------------
# 1. Title (A1) and header-row (row 2) font restyle on both sheets: bold,
#    default size (11), white text - matches the header's dark-blue fill.
# ---------------------------------------------------------------------------
$sheet1Header = $ws1.Range("A2:K2").Font
$sheet1Header.Bold = $true
$sheet1Header.Size = 11
$sheet1Header.Color = 16777215

$sheet1Title = $ws1.Range("A1").Font
$sheet1Title.Bold = $true
$sheet1Title.Size = 11
$sheet1Title.Color = 16777215

$sheet2Header = $ws2.Range("A2:G2").Font
$sheet2Header.Bold = $true
$sheet2Header.Size = 11
$sheet2Header.Color = 16777215

$sheet2Title = $ws2.Range("A1").Font
$sheet2Title.Bold = $true
$sheet2Title.Size = 11
$sheet2Title.Color = 16777215

# ---------------------------------------------------------------------------
# 2. Training Dashboard: refresh PERIOD TO EXPIRE (H) and LAST UPDATE (I)
#    for every data row (3-26). All periods drop by 8 (new check date is
#    8 days later: 16-Sep-2025 instead of 08-Sep-2025).
# ---------------------------------------------------------------------------
$periods = @{
    3  = 470
    4  = 243
    5  = 447
    6  = 240
    7  = 483
    8  = 483
    9  = 219
    10 = 334
    11 = 601
    12 = 338
    13 = 335
    14 = 335
    15 = 260
    16 = 360
    17 = 15
    18 = -103
    19 = -180
    20 = -48
    21 = -48
    22 = 155
    23 = 268
    24 = 313
    25 = 313
    26 = 313
}

foreach ($row in 3..26) {
    $ws1.Cells.Item($row, 8).Value = $periods[$row]
    # Use a quote-prefixed text assignment so the date-like string is kept
    # as literal text instead of being parsed into a date serial value.
    $ws1.Cells.Item($row, 9).Value = "'16-Sep-2025"
}

# ---------------------------------------------------------------------------
# 3. Row 17 (LOTO SOP) is now expired -> flip status and apply the "NOT
#    VALID" look (same red/pink fill + borders used by the other invalid
#    rows, e.g. row 18) by copying its formatting across.
# ---------------------------------------------------------------------------
$ws1.Range("A18:K18").Copy() | Out-Null
$ws1.Range("A17:K17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws1.Range("J17").Value = "NOT VALID"

# ---------------------------------------------------------------------------
# 4. Exam Dashboard: narrow the COMMENTS column and refresh its wording.
# ---------------------------------------------------------------------------
$ws2.Columns.Item(5).ColumnWidth = 14.1666667

foreach ($row in 3..7) {
    $ws2.Cells.Item($row, 5).Value = "date is valid"
}
